$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D updates (row 12-16, 19, 26) -- literal values
$ws.Range("D12").Value = 1398800851.3300049
$ws.Range("D13").Value = 338965688.90999979
$ws.Range("D14").Value = -45752811.059999987
$ws.Range("D16").Value = -50601311.959999993
$ws.Range("D19").Value = -383099999.99999988
$ws.Range("D26").Value = 1009991810.1331247

# Column E updates (row 12-16, 19, 22, 26) -- literal values
$ws.Range("E12").Value = 1361974149
$ws.Range("E13").Value = 332114255
$ws.Range("E14").Value = 6537985
$ws.Range("E15").Value = 2230000000
$ws.Range("E16").Value = 60473972
$ws.Range("E19").Value = 384700000
$ws.Range("E22").Value = 20015625
$ws.Range("E26").Value = 1018613404

# D18 and D21 become formulas (were static values before)
$ws.Range("D18").Formula = "=SUM(D12:D17)"
$ws.Range("D21").Formula = "=SUM(D18:D20)"
